$wb = $excel.ActiveWorkbook

# Rename the "Type" sheet to "ProductType" (sheetId 4 / rId4 unchanged).
$wsType = $wb.Worksheets.Item("Type")
$wsType.Name = "ProductType"

# Make "ProductType" the active sheet (was "Category" before). This flips
# tabSelected from the Category sheetView to the ProductType sheetView and
# updates the workbook's bookViews/workbookView activeTab accordingly.
$wsType.Activate()

Write-Output "done"
